# Insert a new data row at row 6, pushing existing rows 6:122 down to 7:123,
# then populate the newly inserted row 6 with the new weekly reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 6 (Excel copies formatting from
# the row above, which already carries the date number format for column D).
$ws.Rows("6:6").Insert()

$ws.Range("A6").Value = 10
$ws.Range("B6").Value = "Vega Modelo de Temuco"
$ws.Range("C6").Value = "La Araucanía"
$ws.Range("D6").Value = (Get-Date -Year 2023 -Month 9 -Day 21 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E6").Value = 9
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100108
$ws.Range("H6").Value = "Tropicales y subtropicales"
$ws.Range("I6").Value = 100108004
$ws.Range("J6").Value = "Papaya"
$ws.Range("K6").Value = "Cultivar IV Región"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 90
$ws.Range("N6").Value = 24000
$ws.Range("O6").Value = 24000
$ws.Range("P6").Value = 24000
$ws.Range("Q6").Value = "$/bandeja 10 kilos"
$ws.Range("R6").Value = "Provincia del Elquí"
$ws.Range("S6").Value = 2400
$ws.Range("T6").Value = 10
